# Applies odds updates to "Sheet1" for the FlashScore 2024-10-17 workbook.
# Row 9  -> Cusco vs Los Chankas
# Row 13 -> Nacional vs Miramar

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 9 updates ---
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.65
$ws.Range("R9").Value = 2.2

# --- Row 13 updates ---
$ws.Range("G13").Value = 1.36
$ws.Range("H13").Value = 4.33
$ws.Range("I13").Value = 9.5
$ws.Range("L13").Value = 8
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.98
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("W13").Value = 6
$ws.Range("X13").Value = 6
$ws.Range("Z13").Value = 8.5
$ws.Range("AB13").Value = 34
$ws.Range("AD13").Value = 8.5
$ws.Range("AJ13").Value = 26
$ws.Range("AL13").Value = 67
$ws.Range("AM13").Value = 67
$ws.Range("AN13").Value = 3.2
$ws.Range("AO13").Value = 6.5
$ws.Range("AQ13").Value = 19
$ws.Range("AU13").Value = 10
$ws.Range("AW13").Value = 9
$ws.Range("AZ13").Value = 201
